$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (was 45406 -> now 45436, i.e. one month later)
$ws.Range("A1").Value = 45436

# Update prices in D30 and D31
$ws.Range("D30").Value = 1013
$ws.Range("D31").Value = 1917

$wb.Save()
